$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new numeric-looking values would
# otherwise be auto-converted to numbers by Excel (matches original
# inlineStr/text storage in the source file).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated cell values.
$ws.Range("D2").Value = "29.909.04"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "1.634.84"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.94%  "
$ws.Range("D5").Value = "215.40"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("D6").Value = "0.521"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("E7").Value = "  +0.92%  "
$ws.Range("D8").Value = "28.71"
$ws.Range("E8").Value = "  -2.42%  "
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("E10").Value = "  +0.22%  "
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("D12").Value = "1.868.38"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").Value = "1.639.75"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("E14").Value = "  +3.51%  "
$ws.Range("D15").Value = "9.56"
$ws.Range("E15").Value = "  +5.94%  "
$ws.Range("D16").Value = "3.87"
$ws.Range("E16").Value = "  -1.53%  "
$ws.Range("D17").Value = "29.915.92"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "65.51"
$ws.Range("E18").Value = "  +1.75%  "
$ws.Range("D19").Value = "241.70"
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").Value = "0.0₃0705"
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("E21").Value = "  +0.75%  "
$ws.Range("D22").Value = "9.90"
$ws.Range("E22").Value = "  +2.45%  "
$ws.Range("D23").Value = "4.15"
$ws.Range("E23").Value = "  +0.92%  "
$ws.Range("D24").Value = "2.17"
$ws.Range("E24").Value = "  +2.71%  "
$ws.Range("D25").Value = "157.49"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("D26").Value = "15.55"
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("E27").Value = "  -1.32%  "
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.74%  "
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("D31").Value = "1.11"
$ws.Range("E31").Value = "  +1.50%  "
$ws.Range("E32").Value = "  +1.85%  "
$ws.Range("E33").Value = "  -0.56%  "
$ws.Range("D34").Value = "1.424.58"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  +3.23%  "
$ws.Range("E36").Value = "  -1.58%  "
$ws.Range("E37").Value = "  -3.46%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.0171"
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("B39").Value = "HuobiToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D39").Value = "2.29"
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("D40").Value = "75.73"
$ws.Range("E40").Value = "  +8.24%  "
$ws.Range("D41").Value = "0.560"
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("E43").Value = "  +0.78%  "
$ws.Range("E44").Value = "  +1.08%  "
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("D47").Value = "1.775.87"
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("D48").Value = "5.35"
$ws.Range("E48").Value = "  -2.08%  "
$ws.Range("D49").Value = "48.62"
$ws.Range("D50").Value = "93.31"
$ws.Range("E50").Value = "  +5.83%  "
$ws.Range("E51").Value = "  +3.64%  "
